$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: fill in the "Product" (col A) names for the five new FondJoy action
# figures (rows 8-12), in the order they were first typed by the author.
# ---------------------------------------------------------------------------
$ws.Cells.Item(8, 1).Value  = "Thor 1/12 Plastic Model Kit FondJoy"
$ws.Cells.Item(9, 1).Value  = "Spiderman 1/12 Plastic Model Kit FondJoy"
$ws.Cells.Item(10, 1).Value = "Thanos 1/12 Plastic Model Kit FondJoy"
$ws.Cells.Item(12, 1).Value = "Hulk 1/12 Plastic Model Kit FondJoy"
$ws.Cells.Item(11, 1).Value = "Captain America 1/12 Plastic Model Kit FondJoy"

# ---------------------------------------------------------------------------
# Step 2: fill in the "Image" (col D) file names for rows 8-13.
# ---------------------------------------------------------------------------
$ws.Cells.Item(8, 4).Value  = "thor fj.jpg"
$ws.Cells.Item(9, 4).Value  = "spiderman fj.jpg"
$ws.Cells.Item(10, 4).Value = "thanos fj.jpg"
$ws.Cells.Item(11, 4).Value = "cap fj.jpg"
$ws.Cells.Item(12, 4).Value = "hulk fj.jpg"
$ws.Cells.Item(13, 4).Value = "ironman gantry.jpg"

# ---------------------------------------------------------------------------
# Step 3: new "Ironman Gantry" product name (row 13).
# ---------------------------------------------------------------------------
$ws.Cells.Item(13, 1).Value = "Ironman Gantry 1/10 scale"

# ---------------------------------------------------------------------------
# Step 4: a few more image names out of row order, as originally authored.
# ---------------------------------------------------------------------------
$ws.Cells.Item(15, 4).Value = "new superman set.jpg"
$ws.Cells.Item(14, 4).Value = "thunderbolts.jpg"
$ws.Cells.Item(16, 4).Value = "thor love.jpg"

# ---------------------------------------------------------------------------
# Step 5: fix the placeholder image on the small Hogwarts Castle row.
# ---------------------------------------------------------------------------
$ws.Cells.Item(5, 4).Value = "small hog.jpg"

# ---------------------------------------------------------------------------
# Step 6: new Hogwarts Castle (Medium/Large) product names (rows 17-18) and
# their image file names.
# ---------------------------------------------------------------------------
$ws.Cells.Item(17, 1).Value = "Hogwarts Castle (Medium)"
$ws.Cells.Item(18, 1).Value = "Hogwarts Castle (Large)"
$ws.Cells.Item(17, 4).Value = "med hog.jpg"
$ws.Cells.Item(18, 4).Value = "large hog.jpg"

# ---------------------------------------------------------------------------
# Step 7: last three new minifigure-set product names (rows 14-16).
# ---------------------------------------------------------------------------
$ws.Cells.Item(15, 1).Value = "New Superman Movie Set (6 minifigure)"
$ws.Cells.Item(16, 1).Value = "Thor - Love & Thunder Set (8 minifigure)"
$ws.Cells.Item(14, 1).Value = "Thunderbolts Set (8 minifigure)"

# ---------------------------------------------------------------------------
# Step 8: fill in the Price (col B) and Category (col C) values for all the
# newly added rows. These reuse already-existing shared strings / numbers,
# so the order here does not affect the shared string table layout.
# ---------------------------------------------------------------------------
$ws.Cells.Item(8, 2).Value  = 1700
$ws.Cells.Item(8, 3).Value  = "Action Figure"

$ws.Cells.Item(9, 2).Value  = 1700
$ws.Cells.Item(9, 3).Value  = "Action Figure"

$ws.Cells.Item(10, 2).Value = 1700
$ws.Cells.Item(10, 3).Value = "Action Figure"

$ws.Cells.Item(11, 2).Value = 1700
$ws.Cells.Item(11, 3).Value = "Action Figure"

$ws.Cells.Item(12, 2).Value = 1700
$ws.Cells.Item(12, 3).Value = "Action Figure"

$ws.Cells.Item(13, 2).Value = 10500
$ws.Cells.Item(13, 3).Value = "Action Figure"

$ws.Cells.Item(14, 2).Value = 1600
$ws.Cells.Item(14, 3).Value = "Minifigure set"

$ws.Cells.Item(15, 2).Value = 1200
$ws.Cells.Item(15, 3).Value = "Minifigure set"

$ws.Cells.Item(16, 2).Value = 1600
$ws.Cells.Item(16, 3).Value = "Minifigure set"

$ws.Cells.Item(17, 2).Value = 3000
$ws.Cells.Item(17, 3).Value = "Building set"

$ws.Cells.Item(18, 2).Value = 4500
$ws.Cells.Item(18, 3).Value = "Building set"

# ---------------------------------------------------------------------------
# Step 9: column widths widened (best-fit) to fit the new, longer text.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 42.8
$ws.Columns.Item(4).ColumnWidth = 20.0

# ---------------------------------------------------------------------------
# Step 10: move the selection to reflect where the author finished editing.
# ---------------------------------------------------------------------------
$ws.Range("A14").Select()
